# Add two new columns, I (I0) and J (IF), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - copy the existing header formatting (style) from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-24: values for columns I (I0) and J (IF)
$data = @(
    @{ Row = 2;  I0 = 1; IF = 4 },
    @{ Row = 3;  I0 = 1; IF = 6 },
    @{ Row = 4;  I0 = 1; IF = 7 },
    @{ Row = 5;  I0 = 1; IF = 4 },
    @{ Row = 6;  I0 = 1; IF = 5 },
    @{ Row = 7;  I0 = 1; IF = 2 },
    @{ Row = 8;  I0 = 1; IF = 3 },
    @{ Row = 9;  I0 = 8; IF = 8 },
    @{ Row = 10; I0 = 7; IF = 7 },
    @{ Row = 11; I0 = 8; IF = 9 },
    @{ Row = 12; I0 = 2; IF = 3 },
    @{ Row = 13; I0 = 6; IF = 7 },
    @{ Row = 14; I0 = 5; IF = 6 },
    @{ Row = 15; I0 = 7; IF = 7 },
    @{ Row = 16; I0 = 4; IF = 5 },
    @{ Row = 17; I0 = 5; IF = 7 },
    @{ Row = 18; I0 = 7; IF = 8 },
    @{ Row = 19; I0 = 4; IF = 6 },
    @{ Row = 20; I0 = 8; IF = 8 },
    @{ Row = 21; I0 = 2; IF = 4 },
    @{ Row = 22; I0 = 1; IF = 5 },
    @{ Row = 23; I0 = 3; IF = 5 },
    @{ Row = 24; I0 = 5; IF = 6 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I0
    $ws.Cells.Item($r, 10).Value = $entry.IF
}
